{"js": "// Fix contact information missing from short resumes:\n// Insert a new centered paragraph with contact info right after the title\n// paragraph (\"Dheeraj Chand\") and before the \"PROFESSIONAL SUMMARY\" heading.\n// The new paragraph/run must carry NO inherited character formatting (no\n// bold, no explicit size) - just paragraph-level center alignment - so we\n// build the paragraph from clean OOXML rather than splitting the existing\n// (bold, size-28) title run, which would otherwise inherit its rPr.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the title paragraph (first paragraph, holding \"Dheeraj Chand\").\nlet titlePara = null;\nfor (const p of paragraphs.items) {\n  if (p.text.trim() === \"Dheeraj Chand\") {\n    titlePara = p;\n    break;\n  }\n}\n\nif (!titlePara) {\n  throw new Error(\"Could not find the 'Dheeraj Chand' title paragraph\");\n}\n\nconst contactText =\n  \"202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX\";\n\nconst flatOpc = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr>\n              <w:jc w:val=\"center\"/>\n            </w:pPr>\n            <w:r>\n              <w:t>${contactText}</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\n// Insert the new paragraph immediately after the whole title paragraph so it\n// lands as its own sibling block, not merged into / splitting the title run.\nconst titleRange = titlePara.getRange(\"Whole\");\ntitleRange.insertOoxml(flatOpc, \"After\");\n\nawait context.sync();\n", "ps1": "# Fix contact information missing from short resumes:\n# Insert a new centered paragraph with contact info right after the title\n# paragraph (\"Dheeraj Chand\") and before the \"PROFESSIONAL SUMMARY\" heading.\n# The new paragraph/run must carry NO inherited character formatting (no\n# bold, no explicit size) - just paragraph-level center alignment - so we\n# insert clean OOXML at the paragraph-mark boundary rather than typing into\n# / splitting the existing (bold, size-28) title run, which would otherwise\n# hand its rPr down to the new run.\n\n$d = $word.ActiveDocument\n\n# Locate the title paragraph by its text.\n$findRange = $d.Content\n$found = $findRange.Find.Execute(\"Dheeraj Chand\")\nif (-not $found) {\n    throw \"Could not find the 'Dheeraj Chand' title paragraph\"\n}\n\n# $findRange.End now sits right at the paragraph mark that ends the title\n# paragraph (i.e. just before \"PROFESSIONAL SUMMARY\" begins). Build a fresh,\n# zero-width Range there so the insert lands as a new sibling paragraph\n# instead of merging into / inheriting from the title run.\n$insertAt = $findRange.End\n$insertRange = $d.Range($insertAt, $insertAt)\n\n$contactText = \"202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX\"\n\n$xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:pPr><w:jc w:val=\"center\"/></w:pPr><w:r><w:t>' + $contactText + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n$insertRange.InsertXML($xml)\n"}
